## restaurants.xlsx update
## 1) Fill in the audit result that had just come back for the
##    "Статойл Софийская" restaurant (row with blank E/F/G) - date,
##    status and auditor.
## 2) Re-sort the whole audit table by "Дата аудита" (column E) so the
##    freshly dated rows take their place in chronological order and the
##    still-unaudited rows (blank date) drop to the bottom, exactly like
##    the worksheet's existing AutoFilter sort-state (E1:E127) describes.
## 3) Nudge the viewport/selection to where the analyst was last looking
##    and tidy up the row heights of the two "Череповец" rows that moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. fill in the newly-arrived audit data for row 119 ------------------
$ws.Range("E119").Value = 45774
$ws.Range("F119").Value = "Зеленый "
$ws.Range("G119").Value = "Наран"

# --- 2. re-sort the data block by audit date (ascending, blanks last) -----
$sortRange = $ws.Range("A2:H127")
$keyRange  = $ws.Range("E2:E127")
$sortRange.Sort($keyRange, 1, $null, $null, 1, $null, $null, 1)

# the two rows that now hold the larger-font "Череповец" entries need the
# taller row height that travels with that formatting; the rows they
# vacated go back to the default height
$ws.Rows(126).AutoFit()
$ws.Rows(127).AutoFit()
$ws.Rows(111).RowHeight = 15.75
$ws.Rows(112).RowHeight = 15.75

# --- 3. viewport / selection, matching where the user ended up ------------
$excel.ActiveWindow.ScrollRow = 70
$ws.Range("E119").Select()
